$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4; this shifts existing rows 4-57 down to 5-58
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new weekly record
$ws.Cells.Item(4, 1).Value2 = 10
$ws.Cells.Item(4, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(4, 3).Value2 = "La Araucanía"
$ws.Cells.Item(4, 4).Value2 = 44515
$ws.Cells.Item(4, 5).Value2 = 9
$ws.Cells.Item(4, 6).Value2 = "Fruta"
$ws.Cells.Item(4, 7).Value2 = 100101
$ws.Cells.Item(4, 8).Value2 = "Berries"
$ws.Cells.Item(4, 9).Value2 = 100101001
$ws.Cells.Item(4, 10).Value2 = "Arándano (blue)"
$ws.Cells.Item(4, 11).Value2 = "Sin especificar"
$ws.Cells.Item(4, 12).Value2 = "Primera"
$ws.Cells.Item(4, 13).Value2 = 200
$ws.Cells.Item(4, 14).Value2 = 3500
$ws.Cells.Item(4, 15).Value2 = 3500
$ws.Cells.Item(4, 16).Value2 = 3500
$ws.Cells.Item(4, 17).Value2 = "$/kilo"
$ws.Cells.Item(4, 18).Value2 = "Región del Maule"
$ws.Cells.Item(4, 19).Value2 = 3500
$ws.Cells.Item(4, 20).Value2 = 1
